# Atualizei dados da bibi: corrige métricas de retenção anual
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: cohort 2024, period_index 1 -> num_customers goes from 106 to 107
$ws.Range("C36").Value = 107
$ws.Range("E36").Value = 0.05544041450777202

# Row 37: cohort 2025, period_index 0 -> num_customers/cohort_size go from 648 to 654
$ws.Range("C37").Value = 654
$ws.Range("D37").Value = 654
